$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2 used to hold a phone number; replace it with the account's email address.
$ws.Range("A2").Value = "jerinjamesm@gmail.com"

# Turn A2 into a clickable mailto: hyperlink (this also applies the built-in
# "Hyperlink" cell style: underlined, theme-colored font).
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:jerinjamesm@gmail.com")

# B2 keeps its existing password text.
$ws.Range("B2").Value = "Girlboss27*"

# Move the active selection to B2.
$ws.Range("B2").Select()
